$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

# Update PLC live data values per 2025-10-13 13:53:32 snapshot
$ws.Range("B2").Value = 7185
$ws.Range("C3").Value = 165220
$ws.Range("C4").Value = 156166
$ws.Range("C8").Value = 64.98999999999999
